# Updated Masterdata as per 2nd may Data Refresh
# Append two new master-location rows (Arabic postal-code entries for
# codes 10113 / 10114) mirroring the existing pattern used by rows 118-119.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 120 ---------------------------------------------------------
$ws.Range("A120").Value = 10113
$ws.Range("B120").Value = 10113
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = "الرمز البريدي"
$ws.Range("E120").Value = "BNMR"
$ws.Range("F120").Value = "ara"
$ws.Range("G120").Value = $true
$ws.Range("H120").Value = "superadmin"
$ws.Range("I120").Value = "now()"

# --- Row 121 ---------------------------------------------------------
$ws.Range("A121").Value = 10114
$ws.Range("B121").Value = 10114
$ws.Range("C121").Value = 5
$ws.Range("D121").Value = "الرمز البريدي"
$ws.Range("E121").Value = "BNMR"
$ws.Range("F121").Value = "ara"
$ws.Range("G121").Value = $true
$ws.Range("H121").Value = "superadmin"
$ws.Range("I121").Value = "now()"

# Match the saved cursor/selection state recorded in the workbook after
# the data entry (selecting the first fully-empty row below the table).
[void]$ws.Range("A122:XFD1048576").Select()
